$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (Coin, Link, Price, Volume) keep their exact
# string representation instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.740.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.958.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.33%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '380.25'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.53'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.07%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.599'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.41'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.140'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0846'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.61'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.415.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.48'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.950.93'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.959'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.704.01'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.46'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.65'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.66'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.25%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.14%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +20.04%  '

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.170'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.40%  '

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.16'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.88%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.20%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.04'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.74%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.82%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '52.48'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '34.52'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.31%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0436'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.67%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.07'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.79%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.68'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.34%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.39'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.85'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.89%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '124.80'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.98'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.67%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.280'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +18.43%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.07'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.032.66'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.40%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0330'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.88%  '
